# SIMULATION.docx - "Intersection Ahead sign" paragraph:
#   1. Highlight the "9.)Intersection" lead-in (as three separate runs) red.
#   2. Merge the " Ahead sign - " / body / " (CORRECT P200) (WRONG P300)"
#      runs into a single red-highlighted run, dropping the
#      "(CORRECT P200) (WRONG P300)" suffix.

$d = $word.ActiveDocument

# Locate the paragraph that contains the "Intersection Ahead sign" item
# (there's an earlier, unrelated "Intersection signs:" heading paragraph,
# so search paragraph-by-paragraph rather than the whole document).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*9.)Intersection Ahead sign*") {
        $targetPara = $candidate
        break
    }
}

$paraRange = $targetPara.Range

# Pin down the exact "9.)Intersection" span within that paragraph so the
# highlighting below can't accidentally touch the earlier "Intersection"
# that appears in "Intersection signs:".
$locate = $paraRange.Duplicate
$locate.Find.ClearFormatting()
$locate.Find.Execute("9.)Intersection", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$leadStart = $locate.Start
$leadEnd = $locate.End

# Highlight "9" red.
$rNum = $d.Range($leadStart, $leadEnd)
$rNum.Find.ClearFormatting()
$rNum.Find.Replacement.ClearFormatting()
$rNum.Find.Replacement.Highlight = 6
$rNum.Find.Execute("9", $false, $false, $false, $false, $false, $true, 1, $false, "9", 2) | Out-Null

# Highlight ".)" red.
$rPunct = $d.Range($leadStart, $leadEnd)
$rPunct.Find.ClearFormatting()
$rPunct.Find.Replacement.ClearFormatting()
$rPunct.Find.Replacement.Highlight = 6
$rPunct.Find.Execute(".)", $false, $false, $false, $false, $false, $true, 1, $false, ".)", 2) | Out-Null

# Highlight "Intersection" red.
$rWord = $d.Range($leadStart, $leadEnd)
$rWord.Find.ClearFormatting()
$rWord.Find.Replacement.ClearFormatting()
$rWord.Find.Replacement.Highlight = 6
$rWord.Find.Execute("Intersection", $false, $false, $false, $false, $false, $true, 1, $false, "Intersection", 2) | Out-Null

# Merge the trailing three runs into one red-highlighted run, dropping the
# "(CORRECT P200) (WRONG P300)" suffix.
$oldTail = " Ahead sign - The road you are traveling on intersects a highway ahead. Slow down, look to the right and to the left for other traffic, be prepared to stop and yield to crossing traffic. (CORRECT P200) (WRONG P300)"
$newTail = " Ahead sign - The road you are traveling on intersects a highway ahead. Slow down, look to the right and to the left for other traffic, be prepared to stop and yield to crossing traffic."

$rTail = $paraRange.Duplicate
$rTail.Find.ClearFormatting()
$rTail.Find.Replacement.ClearFormatting()
$rTail.Find.Replacement.Highlight = 6
$rTail.Find.Execute($oldTail, $false, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

Write-Host "Done. Paragraph text now: $($targetPara.Range.Text)"
